$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESCOR")

# Delete entire column E (it duplicated the old column F's 1/0 flag); remaining
# columns (old F -> new E) shift left.
$ws.Columns.Item(5).Delete()

# The hidden _xlnm._FilterDatabase name still references the old $F$116 extent;
# repoint it at the new last column (E) now that the sheet is one column narrower.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ESCOR!_FilterDatabase") {
        $n.RefersTo = "=ESCOR!`$A`$1:`$E`$116"
    }
}

# Update the active cell/selection to match the post-edit state.
$ws.Range("E84").Select()

